# Slide 4, Title placeholder: "About Me" -> "Background About Me"
# (rendered as two runs: "Background About " + "Me", matching the
#  OOXML diff which splits the original single run into two runs.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# "About Me" -> first 6 characters are "About ", replace with
# "Background About " so the remaining original run keeps just "Me".
$prefix = $tr.Characters(1, 6)
$prefix.Text = "Background About "
